$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.613.36"
$ws.Range("E2").Value = "  -0.50%  "

$ws.Range("D3").Value = "2.540.71"
$ws.Range("E3").Value = "  -0.05%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'313.35"
$ws.Range("E5").Value = "  +3.10%  "

$ws.Range("D6").Value = "'95.20"
$ws.Range("E6").Value = "  -2.59%  "

$ws.Range("D7").Value = "'0.579"
$ws.Range("E7").Value = "  +0.30%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "'0.540"
$ws.Range("E9").Value = "  -0.93%  "

$ws.Range("D10").Value = "'36.44"
$ws.Range("E10").Value = "  -0.96%  "

$ws.Range("E11").Value = "  -1.47%  "

$ws.Range("D12").Value = "'7.72"
$ws.Range("E12").Value = "  +1.31%  "

$ws.Range("E13").Value = "  -0.98%  "

$ws.Range("D14").Value = "2.926.53"
$ws.Range("E14").Value = "  -0.22%  "

$ws.Range("D15").Value = "'15.69"
$ws.Range("E15").Value = "  +4.17%  "

$ws.Range("D16").Value = "2.527.04"
$ws.Range("E16").Value = "  -1.04%  "

$ws.Range("D17").Value = "'0.867"
$ws.Range("E17").Value = "  +0.06%  "

$ws.Range("D18").Value = "42.637.86"
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("D19").Value = "'13.10"
$ws.Range("E19").Value = "  -3.53%  "

$ws.Range("D20").Value = "'6.69"
$ws.Range("E20").Value = "  +1.49%  "

$ws.Range("D21").Value = "0.0₃0972"
$ws.Range("E21").Value = "  -1.95%  "

$ws.Range("D22").Value = "'71.19"
$ws.Range("E22").Value = "  -1.13%  "

$ws.Range("D23").Value = "'255.16"
$ws.Range("E23").Value = "  +0.10%  "

$ws.Range("D24").Value = "'2.96"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("E25").Value = "  -1.82%  "

$ws.Range("D26").Value = "'27.58"
$ws.Range("E26").Value = "  -1.92%  "

$ws.Range("D27").Value = "'0.988"
$ws.Range("E27").Value = "  -1.01%  "

$ws.Range("D28").Value = "'2.36"
$ws.Range("E28").Value = "  +12.80%  "

$ws.Range("D29").Value = "'39.83"
$ws.Range("E29").Value = "  +5.26%  "

$ws.Range("E30").Value = "  -1.06%  "

$ws.Range("D31").Value = "'5.96"
$ws.Range("E31").Value = "  -2.92%  "

$ws.Range("D32").Value = "'155.61"
$ws.Range("E32").Value = "  -1.45%  "

$ws.Range("D33").Value = "'20.04"
$ws.Range("E33").Value = "  +2.97%  "

$ws.Range("D34").Value = "'3.42"
$ws.Range("E34").Value = "  +3.40%  "

$ws.Range("D35").Value = "'2.15"
$ws.Range("E35").Value = "  +0.91%  "

$ws.Range("D36").Value = "'0.0796"
$ws.Range("E36").Value = "  -0.39%  "

$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("E38").Value = "  -2.84%  "

$ws.Range("D39").Value = "'24.99"
$ws.Range("E39").Value = "  -1.19%  "

$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("D41").Value = "'2.18"
$ws.Range("E41").Value = "  +5.03%  "

$ws.Range("D42").Value = "'3.41"
$ws.Range("E42").Value = "  -0.65%  "

$ws.Range("E43").Value = "  -1.14%  "

$ws.Range("E44").Value = "  -0.33%  "

$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").Value = "2.039.41"
$ws.Range("E46").Value = "  -2.68%  "

$ws.Range("D47").Value = "'85.57"
$ws.Range("E47").Value = "  -3.12%  "

$ws.Range("E48").Value = "  -0.63%  "

$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "'74.94"
$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.782.79"
$ws.Range("E50").Value = "  -0.28%  "

$ws.Range("E51").Value = "  +0.29%  "
